# chore: adapt column header formatting to respective input file names (#7)
# - rename "<field>_old" -> "<field>_FV2310" and "<field>_new" -> "<field>_FV2404"
#   in the header row
# - turn the sheet's used range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (A1:J1 = "_old" suffix, L1:U1 = "_new" suffix,
#        K1 stays "diff") -------------------------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = [string]$cell.Value2 -replace "_old$", "_FV2310"
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = [string]$cell.Value2 -replace "_new$", "_FV2404"
}

# --- 2. Convert the data range into a proper Excel Table ------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
